$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for plain (unstyled) data cells, used to restore the
# default cell style after forcing a numeric-looking string to remain text.
$defaultDataStyle = $ws.Cells.Item(32, 4).Style

$ws.Cells.Item(2, 4).Value = '64.334.47'
$ws.Cells.Item(2, 5).Value = '  +1.09%  '
$ws.Cells.Item(3, 4).Value = '2.632.37'
$ws.Cells.Item(3, 5).Value = '  +0.37%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '''597.08'
$ws.Cells.Item(5, 4).Style = $defaultDataStyle
$ws.Cells.Item(5, 5).Value = '  +0.32%  '
$ws.Cells.Item(6, 4).Value = '''152.49'
$ws.Cells.Item(6, 4).Style = $defaultDataStyle
$ws.Cells.Item(6, 5).Value = '  +1.61%  '
$ws.Cells.Item(7, 5).Value = '  +0.08%  '
$ws.Cells.Item(8, 4).Value = '''0.590'
$ws.Cells.Item(8, 4).Style = $defaultDataStyle
$ws.Cells.Item(9, 5).Value = '  +5.49%  '
$ws.Cells.Item(10, 4).Value = '''5.83'
$ws.Cells.Item(10, 4).Style = $defaultDataStyle
$ws.Cells.Item(10, 5).Value = '  +2.25%  '
$ws.Cells.Item(11, 4).Value = '''0.394'
$ws.Cells.Item(11, 4).Style = $defaultDataStyle
$ws.Cells.Item(11, 5).Value = '  +3.58%  '
$ws.Cells.Item(12, 5).Value = '  +1.14%  '
$ws.Cells.Item(13, 4).Value = '''28.08'
$ws.Cells.Item(13, 4).Style = $defaultDataStyle
$ws.Cells.Item(13, 5).Value = '  +1.51%  '
$ws.Cells.Item(14, 4).Value = '3.104.84'
$ws.Cells.Item(14, 5).Value = '  +0.32%  '
$ws.Cells.Item(15, 4).Value = '''0.0000172'
$ws.Cells.Item(15, 4).Style = $defaultDataStyle
$ws.Cells.Item(15, 5).Value = '  +14.29%  '
$ws.Cells.Item(16, 4).Value = '64.202.69'
$ws.Cells.Item(16, 5).Value = '  +1.16%  '
$ws.Cells.Item(17, 4).Value = '2.605.15'
$ws.Cells.Item(17, 5).Value = '  -1.36%  '
$ws.Cells.Item(18, 5).Value = '  -0.12%  '
$ws.Cells.Item(19, 4).Value = '''4.78'
$ws.Cells.Item(19, 4).Style = $defaultDataStyle
$ws.Cells.Item(19, 5).Value = '  +2.49%  '
$ws.Cells.Item(20, 4).Value = '''350.76'
$ws.Cells.Item(20, 4).Style = $defaultDataStyle
$ws.Cells.Item(20, 5).Value = '  +1.03%  '
$ws.Cells.Item(21, 4).Value = '''7.10'
$ws.Cells.Item(21, 4).Style = $defaultDataStyle
$ws.Cells.Item(21, 5).Value = '  +3.43%  '
$ws.Cells.Item(22, 5).Value = '  +0.26%  '
$ws.Cells.Item(23, 4).Value = '''67.73'
$ws.Cells.Item(23, 4).Style = $defaultDataStyle
$ws.Cells.Item(23, 5).Value = '  +1.94%  '
$ws.Cells.Item(24, 5).Value = '  -2.03%  '
$ws.Cells.Item(25, 4).Value = '''9.24'
$ws.Cells.Item(25, 4).Style = $defaultDataStyle
$ws.Cells.Item(25, 5).Value = '  +0.30%  '
$ws.Cells.Item(26, 4).Value = '''1.67'
$ws.Cells.Item(26, 4).Style = $defaultDataStyle
$ws.Cells.Item(26, 5).Value = '  -0.56%  '
$ws.Cells.Item(27, 5).Value = '  +1.38%  '
$ws.Cells.Item(28, 4).Value = '''554.35'
$ws.Cells.Item(28, 4).Style = $defaultDataStyle
$ws.Cells.Item(28, 5).Value = '  -3.76%  '
$ws.Cells.Item(29, 5).Value = '  -1.07%  '
$ws.Cells.Item(30, 5).Value = '  -0.12%  '
$ws.Cells.Item(31, 4).Value = '0.0₃0911'
$ws.Cells.Item(31, 5).Value = '  +7.85%  '
$ws.Cells.Item(33, 5).Value = '  +5.13%  '
$ws.Cells.Item(34, 4).Value = '''5.54'
$ws.Cells.Item(34, 4).Style = $defaultDataStyle
$ws.Cells.Item(34, 5).Value = '  +5.40%  '
$ws.Cells.Item(35, 5).Value = '  +1.33%  '
$ws.Cells.Item(36, 4).Value = '''0.423'
$ws.Cells.Item(36, 4).Style = $defaultDataStyle
$ws.Cells.Item(36, 5).Value = '  +3.50%  '
$ws.Cells.Item(37, 4).Value = '''166.10'
$ws.Cells.Item(37, 4).Style = $defaultDataStyle
$ws.Cells.Item(37, 5).Value = '  -1.52%  '
$ws.Cells.Item(38, 4).Value = '''20.07'
$ws.Cells.Item(38, 4).Style = $defaultDataStyle
$ws.Cells.Item(38, 5).Value = '  +3.69%  '
$ws.Cells.Item(39, 5).Value = '  +2.84%  '
$ws.Cells.Item(40, 4).Value = '''0.999'
$ws.Cells.Item(40, 4).Style = $defaultDataStyle
$ws.Cells.Item(40, 5).Value = '  -0.07%  '
$ws.Cells.Item(41, 5).Value = '  +0.00%  '
$ws.Cells.Item(42, 4).Value = '''169.23'
$ws.Cells.Item(42, 4).Style = $defaultDataStyle
$ws.Cells.Item(42, 5).Value = '  +0.44%  '
$ws.Cells.Item(43, 5).Value = '  +4.45%  '
$ws.Cells.Item(44, 4).Value = '''23.26'
$ws.Cells.Item(44, 4).Style = $defaultDataStyle
$ws.Cells.Item(44, 5).Value = '  +8.74%  '
$ws.Cells.Item(45, 5).Value = '  -1.88%  '
$ws.Cells.Item(46, 4).Value = '''2.22'
$ws.Cells.Item(46, 4).Style = $defaultDataStyle
$ws.Cells.Item(46, 5).Value = '  +11.60%  '
$ws.Cells.Item(47, 5).Value = '  +1.99%  '
$ws.Cells.Item(48, 5).Value = '  +1.65%  '
$ws.Cells.Item(49, 5).Value = '  +1.24%  '
$ws.Cells.Item(50, 5).Value = '  +0.22%  '
$ws.Cells.Item(51, 4).Value = '0.0₆0232'
$ws.Cells.Item(51, 5).Value = '  +18.26%  '
